$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.006780028343201
$ws.Range("B1").Value = 2.116267204284668
$ws.Range("C1").Value = 6.550631523132324
$ws.Range("D1").Value = 1.750497102737427
$ws.Range("E1").Value = 1.36799144744873
